# garbage collection video link
# 1) Bump the cached "datetimeFigureOut" date field text from 10/29/2023
#    to 10/30/2023 everywhere it is cached (the slide master and every
#    slide layout's Date placeholder).
# 2) Fix the Racket snippet on slide 3: "(1 2 3 4 5)" -> "(list 1 2 3 4 5)".

$p = $ppt.ActivePresentation

$oldDate = "10/29/2023"
$newDate = "10/30/2023"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master's own Date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every custom (slide) layout's Date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Slide 3: "(define mylist (1 2 3 4 5))" -> "(define mylist (list 1 2 3 4 5))"
$slide3 = $p.Slides.Item(3)
$contentShape = $slide3.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange
$oldSnippet = " (1 2 3 4 5))"
$newSnippet = " (list 1 2 3 4 5))"
$idx = $tr.Text.IndexOf($oldSnippet)
if ($idx -ge 0) {
    $sub = $tr.Characters($idx + 1, $oldSnippet.Length)
    $sub.Text = $newSnippet
}
